# Update T20 (column E) appearance counts for a batch of players as part
# of the "full data scraped for extra batting and bowling fields" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_formats_raw")

# Map of row number -> new T20 (column E) value
$updates = @{
    7  = 3   # C J Bowes
    10 = 47  # M S Chapman
    17 = 7   # M J Henry
    20 = 21  # T W M Latham
    21 = 4   # B G Lister
    23 = 38  # A F Milne
    24 = 47  # D J Mitchell
    25 = 63  # J D S Neesham
    29 = 9   # R Ravindra
    33 = 43  # T L Seifert
    35 = 94  # I S Sodhi
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
